$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Copy the header style from L1 (bold, bordered, centered) onto the new M1 header cell,
# then set its text to the new column name "Serviced by ".
$ws.Cells.Item(1, 12).Copy()
$ws.Cells.Item(1, 13).PasteSpecial(-4122)
$ws.Cells.Item(1, 13).Value = "Serviced by "

# Row 8 previously had blank placeholder cells in F8:K8 - the edit fills them
# with the same "nan" placeholder text used elsewhere in the sheet for blanks.
$row8Cols = 6,7,8,9,10,11
foreach ($c in $row8Cols) {
    $ws.Cells.Item(8, $c).Value = "nan"
}
